$d = $word.ActiveDocument

# The document currently has two consecutive "Title"-styled paragraphs near
# the top: the real title (with the inline picture) and a second paragraph
# ("Nulla ante dui, ...") that was also (incorrectly) styled as "Title" and
# bookmarked. The edit turns that second paragraph into a plain paragraph:
# no paragraph style, no run-level font overrides, and no bookmark. Removing
# the bookmark also renumbers every bookmark id that comes after it.

$targetText = "Nulla ante dui, efficitur ut accumsan id, imperdiet ac urna. " +
    "Duis nec eros non ex posuere scelerisque. Duis non dui quam. Vivamus " +
    "pretium pretium lacus sit amet volutpat. In sollicitudin massa " +
    "euismod, consectetur est in, malesuada sem. Pellentesque ullamcorper " +
    "ligula blandit lacinia cursus. Nunc sit amet quam dapibus, blandit " +
    "lacus in, pellentesque lacus. Morbi varius est sapien, vel imperdiet " +
    "turpis varius vitae. Sed laoreet eu magna vel dictum. Nullam eget " +
    "iaculis nisl, et congue orci. Curabitur hendrerit fermentum sapien " +
    "fringilla vehicula. Ut rutrum pretium ligula in accumsan. Donec sed " +
    "facilisis justo."

# Locate the paragraph that still carries the old "Title" styling/bookmark.
$oldPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq $targetText) {
        $oldPara = $cand
        break
    }
}

$followingPara = $d.Paragraphs.Item($oldPara.Index + 1)

# Build a clean, unstyled paragraph (Normal style, no direct rFonts/sz
# overrides) by adding a fresh paragraph at the very end of the document -
# freshly-added paragraphs come out with bare "<w:pPr><w:rPr/></w:pPr>".
$scratch = $d.Paragraphs.Add()
$scratch.Range.Text = $targetText
$cleanFormattedText = $scratch.Range.FormattedText

# Splice that clean paragraph in right before the paragraph that used to
# follow the old Title paragraph ("What is this" / Heading2).
$insertionPoint = $d.Range($followingPara.Range.Start, $followingPara.Range.Start)
$insertionPoint.FormattedText = $cleanFormattedText

# Remove the old Title-styled paragraph (this also removes its bookmark,
# which automatically renumbers every later bookmark id down by one).
$oldPara.Range.Delete()

# Drop the scratch paragraph used only to mint clean formatting.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Delete()
